# Auto-generated Excel COM-interop script to apply the Famfrit_Profits market-data refresh.
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H:N) for specific leve rows
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, matching the scheduled-runner commit diff.

$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 450.14285
$ws.Range("I33").Value = 281.9091
$ws.Range("K33").Value = 281.9091
$ws.Range("M33").Value = -52.90910000000002
$ws.Range("H127").Value = 1197.6923
$ws.Range("I127").Value = 882.8570999999999
$ws.Range("J127").Value = 1565
$ws.Range("K127").Value = 2648.5713
$ws.Range("L127").Value = 4695
$ws.Range("M127").Value = 2311.4287
$ws.Range("N127").Value = -14615
$ws.Range("H132").Value = 5463.7085
$ws.Range("I132").Value = 5869.5
$ws.Range("K132").Value = 17608.5
$ws.Range("M132").Value = -15078.5

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 3332.25
$ws.Range("I2").Value = 2685.5715
$ws.Range("J2").Value = 4237.6
$ws.Range("K2").Value = 2685.5715
$ws.Range("L2").Value = 4237.6
$ws.Range("M2").Value = -2572.5715
$ws.Range("N2").Value = -4463.6
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()
$ws.Range("H94").Value = 10000
$ws.Range("I94").Value = 10000
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 10000
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -9099
$ws.Range("N94").ClearContents()
$ws.Range("H102").Value = 204380.6
$ws.Range("I102").Value = 289543.84
$ws.Range("J102").Value = 5666.3335
$ws.Range("K102").Value = 289543.84
$ws.Range("L102").Value = 5666.3335
$ws.Range("M102").Value = -287921.84
$ws.Range("N102").Value = -8910.333500000001
$ws.Range("H103").Value = 89990
$ws.Range("I103").Value = 89990
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 89990
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = -88818
$ws.Range("N103").ClearContents()
$ws.Range("H116").Value = 3332.25
$ws.Range("I116").Value = 2685.5715
$ws.Range("J116").Value = 4237.6
$ws.Range("K116").Value = 2685.5715
$ws.Range("L116").Value = 4237.6
$ws.Range("M116").Value = -391.5715
$ws.Range("N116").Value = -8825.6
$ws.Range("H122").Value = 25643414
$ws.Range("I122").Value = 2293.4443
$ws.Range("J122").Value = 83335940
$ws.Range("K122").Value = 6880.3329
$ws.Range("L122").Value = 250007820
$ws.Range("M122").Value = -4430.3329
$ws.Range("N122").Value = -250012720
$ws.Range("H132").Value = 30356012
$ws.Range("I132").Value = 12491.615
$ws.Range("K132").Value = 37474.845
$ws.Range("M132").Value = -34944.845

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 3332.25
$ws.Range("I3").Value = 2685.5715
$ws.Range("J3").Value = 4237.6
$ws.Range("K3").Value = 2685.5715
$ws.Range("L3").Value = 4237.6
$ws.Range("M3").Value = -2571.5715
$ws.Range("N3").Value = -4465.6
$ws.Range("H94").Value = 2999.2
$ws.Range("I94").Value = 3000
$ws.Range("J94").Value = 2999
$ws.Range("K94").Value = 3000
$ws.Range("L94").Value = 2999
$ws.Range("M94").Value = -2549
$ws.Range("N94").Value = -3901
$ws.Range("H107").Value = 3565
$ws.Range("I107").Value = 2608.75
$ws.Range("J107").Value = 5477.5
$ws.Range("K107").Value = 2608.75
$ws.Range("L107").Value = 5477.5
$ws.Range("M107").Value = -688.75
$ws.Range("N107").Value = -9317.5
$ws.Range("H112").Value = 71499.5
$ws.Range("J112").Value = 47999
$ws.Range("L112").Value = 47999
$ws.Range("N112").Value = -50953

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("I31").Value = 2320.5
$ws.Range("J31").Value = 22734138
$ws.Range("K31").Value = 2320.5
$ws.Range("L31").Value = 22734138
$ws.Range("M31").Value = -2025.5
$ws.Range("N31").Value = -22734728
$ws.Range("I34").Value = 2320.5
$ws.Range("J34").Value = 22734138
$ws.Range("K34").Value = 2320.5
$ws.Range("L34").Value = 22734138
$ws.Range("M34").Value = -2118.5
$ws.Range("N34").Value = -22734542
$ws.Range("H58").Value = 1956.6
$ws.Range("I58").Value = 1948.9
$ws.Range("J58").Value = 1972
$ws.Range("K58").Value = 1948.9
$ws.Range("L58").Value = 1972
$ws.Range("M58").Value = -1745.9
$ws.Range("N58").Value = -2378
$ws.Range("H107").Value = 1101.5
$ws.Range("J107").Value = 1992
$ws.Range("L107").Value = 1992
$ws.Range("N107").Value = -5832
$ws.Range("H122").Value = 4787169
$ws.Range("I122").Value = 2556
$ws.Range("K122").Value = 7668
$ws.Range("M122").Value = -5218
$ws.Range("H132").Value = 76882.19
$ws.Range("I132").Value = 102569.75
$ws.Range("K132").Value = 307709.25
$ws.Range("M132").Value = -305179.25
$ws.Range("H134").Value = 2136.6562
$ws.Range("I134").Value = 1738.0385
$ws.Range("K134").Value = 5214.1155
$ws.Range("M134").Value = -2679.1155
$ws.Range("H136").Value = 1956.6
$ws.Range("I136").Value = 1948.9
$ws.Range("J136").Value = 1972
$ws.Range("K136").Value = 5846.700000000001
$ws.Range("L136").Value = 5916
$ws.Range("M136").Value = -3296.700000000001
$ws.Range("N136").Value = -11016

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 938
$ws.Range("I3").Value = 938
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 2814
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -2702
$ws.Range("N3").ClearContents()
$ws.Range("H22").Value = 2221.75
$ws.Range("I22").Value = 795.6667
$ws.Range("J22").Value = 6500
$ws.Range("K22").Value = 2387.0001
$ws.Range("L22").Value = 19500
$ws.Range("M22").Value = -2218.0001
$ws.Range("N22").Value = -19838
$ws.Range("H27").Value = 2221.75
$ws.Range("I27").Value = 795.6667
$ws.Range("J27").Value = 6500
$ws.Range("K27").Value = 2387.0001
$ws.Range("L27").Value = 19500
$ws.Range("M27").Value = -2285.0001
$ws.Range("N27").Value = -19704
$ws.Range("H107").Value = 1684.6666
$ws.Range("J107").Value = 1828.625
$ws.Range("L107").Value = 5485.875
$ws.Range("N107").Value = -9325.875
$ws.Range("H114").Value = 499.8
$ws.Range("I114").Value = 166.33333
$ws.Range("K114").Value = 498.99999
$ws.Range("M114").Value = 2755.00001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 687.1111
$ws.Range("J107").Value = 612.5714
$ws.Range("L107").Value = 612.5714
$ws.Range("N107").Value = -4452.5714

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 907.1539
$ws.Range("I16").Value = 753.7273
$ws.Range("K16").Value = 753.7273
$ws.Range("M16").Value = -583.7273
$ws.Range("H81").Value = 89999.5
$ws.Range("I81").Value = 90000
$ws.Range("K81").Value = 90000
$ws.Range("M81").Value = -89002
$ws.Range("H82").Value = 1929.5333
$ws.Range("I82").Value = 1893.7
$ws.Range("J82").Value = 2001.2
$ws.Range("K82").Value = 1893.7
$ws.Range("L82").Value = 2001.2
$ws.Range("M82").Value = -1532.7
$ws.Range("N82").Value = -2723.2
$ws.Range("H84").Value = 89999.5
$ws.Range("I84").Value = 90000
$ws.Range("K84").Value = 270000
$ws.Range("M84").Value = -265008
$ws.Range("H85").Value = 1929.5333
$ws.Range("I85").Value = 1893.7
$ws.Range("J85").Value = 2001.2
$ws.Range("K85").Value = 1893.7
$ws.Range("L85").Value = 2001.2
$ws.Range("M85").Value = -645.7
$ws.Range("N85").Value = -4497.2

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 72144010
$ws.Range("I100").Value = 101001130
$ws.Range("K100").Value = 202002260
$ws.Range("M100").Value = -202001719
$ws.Range("H136").Value = 2631.8333
$ws.Range("J136").Value = 3498.8333
$ws.Range("L136").Value = 10496.4999
$ws.Range("N136").Value = -15596.4999
